$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection(2)
$pt = $ser.Points(16)
Write-Host $pt
